# Ambermoon Advanced workbook update
# "Nearly finalized cave upper level"
#
# Summary of semantic changes applied:
#  1. GlobalVars sheet: var 226 description changed from "unused" to
#     "You opened the treasure room in dwarf mine" (reused, no longer unused).
#  2. GlobalVars sheet: new var 235 appended describing falling through a
#     cave hole.
#  3. MapChanges sheet: updated the description for the "Grandpa's cellar"
#     map change (cave teleport now deactivated at start, plus a note about
#     making back tiles below the later cave black); row height grew to fit
#     the extra line.
#  4. MapChanges sheet: added a new row describing changes to the
#     "Old dwarf mine" map (global var 226 now set when the treasure room
#     is opened).
#  5. Selection/active-sheet bookkeeping: GlobalVars becomes the active
#     sheet/tab (was Chests before), and the MapChanges sheet's remembered
#     selection moves to B10.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1 & 2. GlobalVars sheet
# ---------------------------------------------------------------------
$globalVars = $wb.Worksheets.Item("GlobalVars")

# Global var 226 is no longer unused - it is now set when the treasure
# room in the dwarf mine is opened.
$globalVars.Range("A9").Value = "226: You opened the treasure room in dwarf mine"

# New global var 235.
$globalVars.Range("A18").Value = "235: You fell through a cave hole so that it create a hole below as well"

# ---------------------------------------------------------------------
# 3 & 4. MapChanges sheet
# ---------------------------------------------------------------------
$mapChanges = $wb.Worksheets.Item("MapChanges")

# Updated description for the Grandpa's cellar map change (row 2).
$mapChanges.Range("C2").Value = "Added teleport to renovated house (with condition)`nAdded wind gate teleport (with condition)`nAdded cave teleport (deactivated at start)`nMade back tiles below later cave black`nNPC Karl can create a wind gate there"
$mapChanges.Rows.Item(2).RowHeight = 75

# New row describing the Old dwarf mine map change.
$mapChanges.Range("A8").Value = 438
$mapChanges.Range("B8").Value = "Old dwarf mine"
$mapChanges.Range("C8").Value = "Global var 226 is now set when you open the treasure room"

# Remembered selection on the MapChanges sheet moves to B10.
$null = $mapChanges.Range("B10").Select()

# ---------------------------------------------------------------------
# 5. Active sheet/tab bookkeeping - GlobalVars becomes the active sheet.
# ---------------------------------------------------------------------
$null = $globalVars.Range("A19").Select()
$null = $globalVars.Activate()
